$wb = $excel.ActiveWorkbook

# Sheet "Users": A2 contact name changed
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Julie Carthane"

# Sheet "GiftLog": B2 contact name changed, C2 gift value changed
$wsGiftLog = $wb.Worksheets.Item("GiftLog")
$wsGiftLog.Range("B2").Value = "Julie Carthane"
$wsGiftLog.Range("C2").Value = "180"

# Restore selections / active cells to match the saved view state
$wsUsers.Range("C20").Select() | Out-Null
$wsGiftLog.Range("O5").Select() | Out-Null

# Make GiftLog the active sheet (tabSelected) as in the source file
$wsGiftLog.Activate() | Out-Null
